$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 136, pushing all existing data (rows 136-256)
# down to rows 137-257.
$ws.Rows(136).Insert()

# Populate the newly inserted row 136 with the new record.
$ws.Cells.Item(136, 1).Value = 5
$ws.Cells.Item(136, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(136, 3).Value = "Maule"
$ws.Cells.Item(136, 4).Value = 44669
$ws.Cells.Item(136, 5).Value = 7
$ws.Cells.Item(136, 6).Value = 100112006
$ws.Cells.Item(136, 7).Value = "Repollo"
$ws.Cells.Item(136, 8).Value = "Crespo record"
$ws.Cells.Item(136, 9).Value = "Primera"
$ws.Cells.Item(136, 10).Value = 2000
$ws.Cells.Item(136, 11).Value = 1000
$ws.Cells.Item(136, 12).Value = 1000
$ws.Cells.Item(136, 13).Value = 1000
$ws.Cells.Item(136, 14).Value = '$/unidad'
$ws.Cells.Item(136, 15).Value = "Región del Maule"
$ws.Cells.Item(136, 16).Value = 1000
$ws.Cells.Item(136, 17).Value = 1
$ws.Cells.Item(136, 18).Value = "Hortaliza"
